# Insert a new data row before row 151 (shifts existing rows 151-221 down
# to 152-222) and populate the new row with the latest "Brocoli" price
# observation, per the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(151).Insert()

$ws.Range("A151").Value = 5
$ws.Range("B151").Value = "Macroferia Regional de Talca"
$ws.Range("C151").Value = "Maule"
$ws.Range("D151").Value = 44510
$ws.Range("E151").Value = 7
$ws.Range("F151").Value = 100112023
$ws.Range("G151").Value = "Brócoli"
$ws.Range("H151").Value = "Sin especificar"
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 4000
$ws.Range("K151").Value = 500
$ws.Range("L151").Value = 500
$ws.Range("M151").Value = 500
$ws.Range("N151").Value = "$/unidad"
$ws.Range("O151").Value = "Región del Maule"
$ws.Range("P151").Value = 500
$ws.Range("Q151").Value = 1
$ws.Range("R151").Value = "Hortaliza"
